$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.105.94"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "1.790.99"
$ws.Range("E3").Value = "  +0.24%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'226.97"
$ws.Range("E5").Value = "  +1.22%  "

$ws.Range("E6").Value = "  -0.54%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "'32.18"
$ws.Range("E8").Value = "  -1.94%  "

$ws.Range("E9").Value = "  +3.11%  "

$ws.Range("D10").Value = "'0.0694"
$ws.Range("E10").Value = "  -2.18%  "

$ws.Range("E11").Value = "  +1.04%  "

$ws.Range("D12").Value = "2.047.25"
$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("D13").Value = "'11.54"
$ws.Range("E13").Value = "  +6.23%  "

$ws.Range("D14").Value = "1.792.70"
$ws.Range("E14").Value = "  +0.22%  "

$ws.Range("D15").Value = "'0.624"
$ws.Range("E15").Value = "  -0.02%  "

$ws.Range("D16").Value = "34.094.68"
$ws.Range("E16").Value = "  +0.14%  "

$ws.Range("D17").Value = "'4.19"
$ws.Range("E17").Value = "  +0.79%  "

$ws.Range("D18").Value = "'68.04"
$ws.Range("E18").Value = "  +0.28%  "

$ws.Range("D19").Value = "'244.73"
$ws.Range("E19").Value = "  -0.23%  "

$ws.Range("D20").Value = "0.0₃0783"
$ws.Range("E20").Value = "  -0.81%  "

$ws.Range("D21").Value = "'10.93"
$ws.Range("E21").Value = "  +1.03%  "

$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").Value = "'4.11"
$ws.Range("E23").Value = "  +0.35%  "

$ws.Range("E24").Value = "  -3.06%  "

$ws.Range("D25").Value = "'162.65"
$ws.Range("E25").Value = "  +1.28%  "

$ws.Range("E26").Value = "  +1.64%  "

$ws.Range("D27").Value = "'16.31"
$ws.Range("E27").Value = "  -0.22%  "

$ws.Range("E28").Value = "  +1.30%  "

$ws.Range("E29").Value = "  +0.23%  "

$ws.Range("E30").Value = "  +1.11%  "

$ws.Range("E31").Value = "  +1.19%  "

$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("E33").Value = "  +3.02%  "

$ws.Range("E34").Value = "  +1.04%  "

$ws.Range("D35").Value = "1.416.98"
$ws.Range("E35").Value = "  +1.46%  "

$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("E37").Value = "  +2.49%  "

$ws.Range("E38").Value = "  -1.01%  "

$ws.Range("E39").Value = "  +5.88%  "

$ws.Range("D40").Value = "'80.94"
$ws.Range("E40").Value = "  +3.64%  "

$ws.Range("D42").Value = "'0.920"
$ws.Range("E42").Value = "  +0.49%  "

$ws.Range("E43").Value = "  -0.17%  "

$ws.Range("D44").Value = "'13.43"
$ws.Range("E44").Value = "  +7.04%  "

$ws.Range("E45").Value = "  +3.26%  "

$ws.Range("E46").Value = "  +1.97%  "

$ws.Range("E47").Value = "  -0.48%  "

$ws.Range("E48").Value = "  -5.69%  "

$ws.Range("D49").Value = "'107.42"
$ws.Range("E49").Value = "  -0.66%  "

$ws.Range("D50").Value = "1.948.78"
$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("E51").Value = "  +0.13%  "
